$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 10.934894
$ws.Range("H2").Value = 32.804682
$ws.Range("I2").Value = 0.3698068269583527
$ws.Range("J2").Value = 0.3698068269583527
$ws.Range("Q2").Value = 16.04703713422266
$ws.Range("R2").Value = 144.423334208004
$ws.Range("S2").Value = 0.06992139188014831
$ws.Range("T2").Value = 0.06992139188014831

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 10.934894
$ws.Range("H3").Value = 32.804682
$ws.Range("I3").Value = 0.3698068269583527
$ws.Range("J3").Value = 0.3698068269583527
$ws.Range("Q3").Value = 68.82404058776666
$ws.Range("R3").Value = 619.4163652899
$ws.Range("S3").Value = 0.2998854350782044
$ws.Range("T3").Value = 0.2998854350782044

# Row 4
$ws.Range("I4").Value = 0.3872921463699351
$ws.Range("J4").Value = 0.3872921463699351
$ws.Range("Q4").Value = 16.805778048254
$ws.Range("S4").Value = 0.07322743650020749
$ws.Range("T4").Value = 0.07322743650020749

# Row 5
$ws.Range("I5").Value = 0.3872921463699351
$ws.Range("J5").Value = 0.3872921463699351
$ws.Range("Q5").Value = 72.07819990865001
$ws.Range("R5").Value = 648.7037991778501
$ws.Range("S5").Value = 0.3140647098697277
$ws.Range("T5").Value = 0.3140647098697276

# Row 6
$ws.Range("I6").Value = 0.2429010266717122
$ws.Range("J6").Value = 0.2429010266717122
$ws.Range("S6").Value = 0.04592662070004419
$ws.Range("T6").Value = 0.04592662070004419

# Row 7
$ws.Range("I7").Value = 0.2429010266717122
$ws.Range("J7").Value = 0.2429010266717122
$ws.Range("S7").Value = 0.196974405971668
$ws.Range("T7").Value = 0.1969744059716679
